# Daily attendance processing - 2025-12-25 06:43:20
# Normalizes the "Recorded By" (column G) values so that the system
# account name is listed after the other recorder(s) rather than first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, system, backup@backdoor.com") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}
